$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data to append below the existing data (rows 6-10)
$rows = @(
    @("Charles Ferreira", 0.78, 0.5600000000000001, 4.24, 0.77, 0.8, 1.43),
    @("testando", 2.42, 5, 4.95, 1.75, 2.56, 3.34),
    @("b", 0.91, 0.68, 1.05, 1.01, 0.8100000000000001, 0.89),
    @("a", 0.44, 1.1, 0.63, 0.76, 0.8100000000000001, 0.75),
    @("deepseek", 4.5, 5, 4.75, 3.39, 1.78, 3.88)
)

$startRow = 6
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    for ($c = 1; $c -le $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowData[$c - 1]
    }
}
